# "Generate Report for Handback"
#
# This script updates the zh-cn / de-de localization-status worksheets (and
# the roll-up Overview sheet) to reflect that both files have been handed
# back and are in sync with en-US:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - "Latest Target File" is now populated with a hyperlink to the source .md
#   - "Latest Handback File" is populated with the xlf file used for handback
#   - "Latest Handback DateTime" is stamped with the handback time

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab58b6a29793ea947d241308fbd4b42aef3623f1/e2e/6bc5fd6b-83dd-44ab-a817-8de84405dc24.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab58b6a29793ea947d241308fbd4b42aef3623f1/e2e/9f75aa24-c862-4956-be99-e0355a2c7a3b.md"

$md1 = "6bc5fd6b-83dd-44ab-a817-8de84405dc24.md"
$md2 = "9f75aa24-c862-4956-be99-e0355a2c7a3b.md"

# ---------------------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status columns for both rows
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value2 = $statusText
$wsOverview.Range("F2").Value2 = $statusText
$wsOverview.Range("E3").Value2 = $statusText
$wsOverview.Range("F3").Value2 = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.14437139602
$wsOverview.Columns.Item(6).ColumnWidth = 29.14437139602

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("C2").Value2 = $statusText
$wsZhCn.Range("C3").Value2 = $statusText

$zhXlf1 = $wsZhCn.Range("G2").Value2
$zhXlf2 = $wsZhCn.Range("G3").Value2

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $urlMd1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $md1) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlMd2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $md2) | Out-Null

$wsZhCn.Range("J2").Value2 = $zhXlf1
$wsZhCn.Range("J3").Value2 = $zhXlf2

$wsZhCn.Range("K2").Value2 = "2016-09-03 03:05:54"
$wsZhCn.Range("K3").Value2 = "2016-09-03 03:05:54"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.14437139602
$wsZhCn.Columns.Item(9).ColumnWidth = 39.16666666667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.16666666667

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe.Range("C2").Value2 = $statusText
$wsDeDe.Range("C3").Value2 = $statusText

$deXlf1 = $wsDeDe.Range("G2").Value2
$deXlf2 = $wsDeDe.Range("G3").Value2

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $urlMd1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $md1) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlMd2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $md2) | Out-Null

$wsDeDe.Range("J2").Value2 = $deXlf1
$wsDeDe.Range("J3").Value2 = $deXlf2

$wsDeDe.Range("K2").Value2 = "2016-09-03 03:06:04"
$wsDeDe.Range("K3").Value2 = "2016-09-03 03:06:04"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.14437139602
$wsDeDe.Columns.Item(9).ColumnWidth = 39.16666666667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.16666666667
